$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "ItemData.Item.1000004"
$ws.Range("B10").Value = "솔라리"
$ws.Range("C10").Value = ""
$ws.Range("C10").Font.Bold = $false
$ws.Range("D10").Value = ""
$ws.Range("D10").Font.Bold = $false

$ws.Range("A11").Value = "MapNpcData.MapNpcMenu.1000005.1"
$ws.Range("B11").Value = "컷신1 재성"
$ws.Range("C11").Value = ""
$ws.Range("C11").Font.Bold = $false
$ws.Range("D11").Value = ""
$ws.Range("D11").Font.Bold = $false

$ws.Range("A12").Value = "MapNpcData.MapNpcMenu.1000005.2"
$ws.Range("B12").Value = "컷신2 재생"
$ws.Range("C12").Value = ""
$ws.Range("C12").Font.Bold = $false
$ws.Range("D12").Value = ""
$ws.Range("D12").Font.Bold = $false
